# Snake_River_TRT_Populations.xlsx - fix a few data errors
#
# 1. Row 8  (Steelhead/Salmon): TRT code "MFPAN-s" was wrong, should be "SRPAN-s"
# 2. Row 25 (Steelhead/Hells Canyon, SNHCT-s): Status should be "Extirpated", not "Extant"
# 3. Row 41 (Chinook/Middle Fork Salmon, MFPIS "Pistol Creek"): this row is a duplicate/erroneous
#    entry and is removed entirely
# 4. Row 46 (Chinook, SRCHA "Chamberlain Creek"): MPG_DPS should be "Middle Fork Salmon",
#    not its own "Chamberlain Creek" group
# 5. Row 50 (Chinook, SRLSR "Little Salmon River"): MPG_DPS should be "South Fork Salmon",
#    not "Little Salmon"
# 6. Row 57 (Chinook/Grande Ronde Imnaha, IRBSH "Big Sheep Creek"): Status should be
#    "Extirpated", not "Extant"
# 7. Row 60 (Chinook/Lower Snake, SNASO "Asotin Creek"): Status should be "Extirpated",
#    not "Extant"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix TRT code for Panther Creek steelhead population
$ws.Range("C8").Value = "SRPAN-s"

# 2. Hells Canyon tributaries steelhead population is extirpated
$ws.Range("E25").Value = "Extirpated"

# 3. Remove the erroneous Pistol Creek (MFPIS) row entirely; everything below shifts up
$ws.Rows.Item(41).Delete()

# After the deletion above, the rows that used to be 46/50/57/60 are now 45/49/56/59
# 4. Chamberlain Creek Chinook population belongs to the Middle Fork Salmon MPG
$ws.Range("B45").Value = "Middle Fork Salmon"

# 5. Little Salmon River Chinook population belongs to the South Fork Salmon MPG
$ws.Range("B49").Value = "South Fork Salmon"

# 6. Big Sheep Creek Chinook population is extirpated
$ws.Range("E56").Value = "Extirpated"

# 7. Asotin Creek Chinook population is extirpated
$ws.Range("E59").Value = "Extirpated"

# Leave the selection where the editor ended up after making these fixes
$ws.Range("D43").Select()
